# Update the "Estado de Cuenta" detail table (rows 16-29):
# previous debt periods are removed and replaced with a new set of periods,
# per the commit "Elimna EC anteriores y se agregan nuevos, se modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each row: TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$rows = @(
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2108", 116000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2107", 120000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2106", 120000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2105", 120000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2104", 120000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2103", 120000, 3000000),
    @("CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2102", 120000, 3000000),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2108", 33942, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2107", 35112, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2106", 35112, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2105", 35112, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2104", 35112, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2103", 35112, 908526),
    @("CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2102", 35112, 908526)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G - Salario Basico
}
